# Update crypto symbol list (coinranking.com scrape refresh)
# Commit: "Updated symbol list on Mon Feb 13 17:28:43 UTC 2023 with GitHub Actions"
#
# The worksheet stores every value (prices, percent changes, URLs, coin names)
# as text, even when it looks numeric. Excel would otherwise auto-convert
# numeric-looking text (e.g. "285.32" or "-10.33%") into real numbers/percentages
# on assignment and silently drop significant trailing zeros (e.g. "0.1050").
# To prevent that, every cell is forced to Text format ("@") before the value is
# assigned, and then its Style is reset to "Normal" so no stray formatting is left
# behind (matches the un-styled cells in the original workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is @("<CellAddress>", "<NewValue>")
$cellUpdates = @(
    @("D2", '285.32'),
    @("E2", '-10.33%'),
    @("D3", '39.94'),
    @("E3", '-3.68%'),
    @("E4", '-4.10%'),
    @("D5", '0.07270'),
    @("E5", '-6.11%'),
    @("B6", 'FTXToken'),
    @("C6", 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'),
    @("D6", '1.510'),
    @("E6", '-10.77%'),
    @("B7", 'MXToken'),
    @("C7", 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @("D7", '0.9137'),
    @("E7", '-4.01%'),
    @("B8", 'LiechtensteinCryptoassetsExchange'),
    @("C8", 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @("D8", '0.1198'),
    @("E8", '-5.39%'),
    @("B9", 'WazirX'),
    @("C9", 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @("D9", '0.1704'),
    @("E9", '-7.10%'),
    @("B10", 'MandalaExchangeToken'),
    @("C10", 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @("D10", '0.08631'),
    @("E10", '-5.72%'),
    @("B11", 'BitrueCoin'),
    @("C11", 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @("D11", '0.04175'),
    @("E11", '-5.29%'),
    @("B12", 'BitMartToken'),
    @("C12", 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @("D12", '0.1050'),
    @("E12", '-0.01%'),
    @("B13", 'BitForexToken'),
    @("C13", 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @("D13", '0.001268'),
    @("E13", '-1.28%'),
    @("B14", 'TigerCash'),
    @("C14", 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @("D14", '0.005955'),
    @("E14", '-1.19%'),
    @("B15", 'LEO'),
    @("C15", 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @("D15", '3.409'),
    @("E15", '2.15%'),
    @("B16", 'GateToken'),
    @("C16", 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @("D16", '4.302'),
    @("E16", '-0.51%'),
    @("E18", '-2.08%'),
    @("D19", '7.857'),
    @("E19", '2.55%'),
    @("D20", '0.1342'),
    @("E20", '-0.72%'),
    @("D21", '0.2885'),
    @("E21", '2.28%'),
    @("D22", '0.03842'),
    @("E22", '-4.51%'),
    @("E23", '0.39%'),
    @("D24", '0.003790'),
    @("D25", '0.0001285'),
    @("E25", '1.22%'),
    @("D26", '0.0003729'),
    @("D38", '0.02301'),
    @("E38", '-9.65%'),
    @("D39", '0.04928'),
    @("E39", '-7.97%'),
    @("D40", '0.007420'),
    @("E40", '272.87%'),
    @("D41", '0.007696'),
    @("E41", '-1.05%'),
    @("D42", '0.1264'),
    @("E42", '-4.12%'),
    @("D43", '0.007392'),
    @("E43", '0.67%'),
    @("D44", '0.006941'),
    @("E44", '-8.27%'),
    @("D45", '0.3089'),
    @("E45", '-10.04%'),
    @("D46", '0.00006375'),
    @("E46", '-4.68%'),
    @("D47", '0.00000000753'),
    @("E47", '0.43%'),
    @("E48", '30.28%'),
    @("E49", '0.09%'),
    @("D50", '0.00002108'),
    @("E50", '0.43%'),
    @("D51", '0.0002008'),
    @("E51", '0.43%')
)

foreach ($update in $cellUpdates) {
    $address = $update[0]
    $newValue = $update[1]
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}
